$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the Version value from 1.0.0 to 1.1.0
$ws.Range("B3").Value = "1.1.0"

# Update the Date value
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# Re-assert the wrap-text alignment on the header/body styles so the
# "applyAlignment" bit is explicitly persisted (it was present before only
# implicitly, i.e. WrapText read back as False until this is (re)applied).
$ws.Range("A1:B1").WrapText = $true
$ws.Range("A2:B14").WrapText = $true

$ws2 = $wb.Worksheets.Item("Include from FSIII")
$ws2.Range("A1:C1").WrapText = $true
$ws2.Range("A2:C2").WrapText = $true
$ws2.Range("A3:B3").WrapText = $true
$ws2.Range("A4:B4").WrapText = $true
